{"js": "// Merge the word-by-word split runs in the Title, Author and Abstract\n// paragraphs into a single run each, with the same combined text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// Map of paragraph style -> exact replacement text (matches the\n// concatenation of the existing split runs, so visible content is\n// unchanged; only the run segmentation collapses to one run).\nconst replacements = {\n  \"Title\": \"Questions: Logarithms\",\n  \"Author\": \"Zo\u00eb Gemmell\",\n  \"Abstract\": \"A selection of questions for the study guide on logarithms.\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const style = para.style;\n  if (Object.prototype.hasOwnProperty.call(replacements, style)) {\n    para.insertText(replacements[style], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Merge the word-by-word split runs in the Title, Author and Abstract\n# paragraphs into a single run each (same combined text, just no\n# longer split run-per-word).\n$doc = $word.ActiveDocument\n\n$targets = @{\n    \"Title\"    = \"Questions: Logarithms\"\n    \"Author\"   = \"Zo\u00eb Gemmell\"\n    \"Abstract\" = \"A selection of questions for the study guide on logarithms.\"\n}\n\nforeach ($p in $doc.Paragraphs) {\n    $style = $p.Range.Style.NameLocal\n    if ($targets.ContainsKey($style)) {\n        $newText = $targets[$style]\n        $r = $p.Range\n        # Paragraph range text includes the trailing paragraph mark;\n        # drop it so Find matches exactly the visible run text.\n        $oldText = $r.Text.TrimEnd([char]13)\n        $r.Find.ClearFormatting()\n        $r.Find.Replacement.ClearFormatting()\n        $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    }\n}\n"}
